# Aula 01 Teoria de Compiladores - slide "Exemplos: ..." paragraph
# Change: "...lavadora de roupas, etc" -> "...lavadora de roupas, etc."
# (the trailing comma moves from the end of the preceding run to the
# start of the final run, and a period is appended after "etc")

$p = $ppt.ActivePresentation

$targetShape = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shp = $slide.Shapes.Item($shi)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -like "*lavadora de roupas*") {
                $targetShape = $shp
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange

# Select ", etc" (the comma+space that trails "...lavadora de roupas"
# together with the word "etc" that follows) as a single range -- this
# span straddles the original run boundary, so re-assigning its text
# splits the paragraph into a run ending at "...roupas" (unchanged) and
# a freshly created run holding the replacement text ", etc.".
$fullText = $tr.Text
$spanIdx = $fullText.IndexOf(", etc")
$spanRange = $tr.Characters($spanIdx + 1, 5)
$spanRange.Text = ", etc."
